$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Borrado logico de proveedores:
# Row 31 (proveedores pagos - related item) status changes from "en proceso" to "terminado"
$ws.Range("B31").Value = "terminado"

# New task row appended at row 36
$ws.Range("A36").Value = "importar tabla articulos con precios y costos en moneda adecuada"
$ws.Range("B36").Value = "no comenzado"

# Update selection to match the diff (A40)
$ws.Range("A40").Select()
